$d = $word.ActiveDocument

# --- 1. Insert a new "Meta description" paragraph right after the title (Heading1) ---
$titlePara = $d.Paragraphs.Item(1)
$insertPoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Beat the Beast Mighty Sphinx and play for free. Enjoy high volatility, free spins, and traditional gameplay.</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$insertPoint.InsertXML($metaXml)

# Remove the stray blank paragraph that InsertXML leaves behind as a separator
$blankPara = $d.Paragraphs.Item(3)
[void]$blankPara.Range.Delete()

# --- 2. Remove the duplicated title paragraph (bold) near the end of the document ---
$dupTitleIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Play Beat the Beast Mighty Sphinx for Free: Review`r") {
        $dupTitleIndex = $i
    }
}
$dupTitlePara = $d.Paragraphs.Item($dupTitleIndex)
[void]$dupTitlePara.Range.Delete()

# --- 3. Replace the text of the final (italic) paragraph with the new image prompt ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$lastRange.Text = 'Prompt: DALLE, create a cartoon-style feature image for "Beat the Beast Mighty Sphinx" featuring a happy Maya warrior with glasses. Description: The feature image should be in cartoon-style with bright and vibrant colors. The main element of the image should be a happy and confident Maya warrior with glasses, standing in front of a giant Sphinx. The warrior should be wearing traditional Maya clothing, with a feather headdress and accessories. The background should have an Egyptian theme, with hieroglyphics and pyramids visible. The Sphinx should be portrayed as dark and imposing, with glowing yellow eyes. The image should emphasize the adventure, excitement, and mystery of Ancient Egypt, while also showcasing the unique blend of Maya and Egyptian elements in the game.'
